$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 14) to the bottom of the portfolio table.
# The Date column needs to stay plain text (matching the existing rows),
# so force a text number format while assigning the value, then restore
# the cell's style back to the default "Normal" so no extra formatting
# is left behind on the cell.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2025-08-29"
$ws.Range("A14").Style = "Normal"

$ws.Range("B14").Value = 56.43000030517578
$ws.Range("C14").Value = 669
$ws.Range("D14").Value = 313.9500122070312
